# Update countries & provincias Spain
# Applies the 27-May-2020 20:05 -> 20:35 data refresh to the "Pais" sheet:
#  - swap the Uganda / Yemen rows (their labels traded places in the
#    shared-string table) and refresh their stats
#  - refresh case/death counters for Estados Unidos, Espana, Alemania,
#    India, Ecuador, Suiza, Oman and Marruecos
#  - bump the "Datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp string -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 20:35"

# --- swap Uganda (row 153) and Yemen (row 154) labels ------------------
$ws.Range("A153").Value = "Yemen"
$ws.Range("A154").Value = "Uganda"

# --- refresh numeric stats ----------------------------------------------
# columns: B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

# row 4 - Estados Unidos
$ws.Range("B4").Value = 1736743
$ws.Range("C4").Value = 11468
$ws.Range("D4").Value = 483090
$ws.Range("E4").Value = 1152183
$ws.Range("G4").Value = 898
$ws.Range("H4").Value = 101470

# row 7 - Espana
$ws.Range("B7").Value = 283849
$ws.Range("C7").Value = 510
$ws.Range("E7").Value = 59773
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 27118

# row 11 - Alemania
$ws.Range("B11").Value = 181757
$ws.Range("C11").Value = 469
$ws.Range("E11").Value = 10434
$ws.Range("G11").Value = 25
$ws.Range("H11").Value = 8523

# row 13 - India
$ws.Range("B13").Value = 158077
$ws.Range("C13").Value = 7284
$ws.Range("D13").Value = 67749
$ws.Range("E13").Value = 85794

# row 27 - Ecuador
$ws.Range("B27").Value = 38103
$ws.Range("C27").Value = 748
$ws.Range("D27").Value = 18425
$ws.Range("E27").Value = 16403
$ws.Range("G27").Value = 72
$ws.Range("H27").Value = 3275

# row 32 - Suiza
$ws.Range("D32").Value = 28300
$ws.Range("E32").Value = 559

# row 58 - Oman
$ws.Range("E58").Value = 6157
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 39

# row 62 - Marruecos
$ws.Range("B62").Value = 7601
$ws.Range("C62").Value = 24
$ws.Range("D62").Value = 4978
$ws.Range("E62").Value = 2421

# row 153 - now Yemen
$ws.Range("B153").Value = 256
$ws.Range("C153").Value = 7
$ws.Range("D153").Value = 10
$ws.Range("E153").Value = 193
$ws.Range("G153").Value = 4
$ws.Range("H153").Value = 53

# row 154 - now Uganda
$ws.Range("B154").Value = 253
$ws.Range("D154").Value = 69
$ws.Range("E154").Value = 184
$ws.Range("H154").Value = 0
